# Homework3.docx edit: continued work on hw3
# Applies the targeted paragraph replacements described by the diff:
#  - Para 17 (Problem 1 list item 5): split "stopping criteria" -> "stopping criterion" across runs
#  - Para 24 (Problem 2 list item 3): TODO -> Bagging/Boosting strengths text
#  - Para 30 (Problem 3 list item 3): add <w:lastRenderedPageBreak/>
#  - Para 31 (Problem 3 list item 4): remove <w:lastRenderedPageBreak/>, split "dropout, and" -> "dropout and" across runs
#  - Para 36 (Problem 4 list item 1): TODO -> Dimensionality reduction text
#  - Para 43 (Problem 5 list item 1): TODO -> Output/Hidden/Input layer activations text

$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $xml) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.InsertXML($xml)
}

$xml17 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1CF197A7" w14:textId="671505EC" w:rsidR="0073091B" w:rsidRDefault="0073091B" w:rsidP="00B270C8"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:t>possible stopping criterion</w:t></w:r><w:r><w:t xml:space="preserve"> for this process can be when there are no more input features to split on since there are only 3 options. </w:t></w:r></w:p>'
Replace-ParagraphXml 17 $xml17

$xml24 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4DB1D5C2" w14:textId="39D2D97D" w:rsidR="0073091B" w:rsidRDefault="00025C4A" w:rsidP="0073091B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="15"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Some strengths of Bagging is that the model is robust against outliers, less likely to overfit, and does not require </w:t></w:r><w:r><w:t>advanced parameter tuning. However, Bagging can be computationally expensive. Boosting puts more weight on weak classifiers from the previous phase so it is more vulnerable to overfitting and outliers. Both of these ensemble methods are less interpretable than the original models they are built upon.</w:t></w:r></w:p>'
Replace-ParagraphXml 24 $xml24

$xml30 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0590DB86" w14:textId="7A4B2E81" w:rsidR="00025C4A" w:rsidRDefault="00561040" w:rsidP="0073091B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">No, increasing the number of features in the model would make the model more complex which would make it more likely to overfit. </w:t></w:r></w:p>'
Replace-ParagraphXml 30 $xml30

$xml31 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3D7DDAE3" w14:textId="4FC267C5" w:rsidR="00025C4A" w:rsidRDefault="00025C4A" w:rsidP="0073091B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Randomly zeroing out half the nodes in the network would help with overfitting. This technique is known as </w:t></w:r><w:r><w:t>dropout and</w:t></w:r><w:r><w:t xml:space="preserve"> is equivalent to training different neural networks and averaging their effects. </w:t></w:r></w:p>'
Replace-ParagraphXml 31 $xml31

$xml36 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7E867A61" w14:textId="1AFBBCA1" w:rsidR="00025C4A" w:rsidRDefault="00025C4A" w:rsidP="00025C4A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Dimensionality reduction can help with overfitting since the model has fewer degrees of freedom. Dimensionality reduction can also help to reduce data storage space, computation time, and remove redundant features. Lastly, dimensionality reduction can be useful for visualizations of high dimensionality data. After dimensionality reduction, the input features are less interpretable because the new axis may not correspond to a real world measurement. </w:t></w:r></w:p>'
Replace-ParagraphXml 36 $xml36

$xml43 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1E7881C8" w14:textId="46842328" w:rsidR="00561040" w:rsidRDefault="00151DEF" w:rsidP="00151DEF"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Output </w:t></w:r><w:r><w:t xml:space="preserve">layers can use </w:t></w:r><w:r><w:t xml:space="preserve">sigmoid (binary), softmax (multiclass), </w:t></w:r><w:r><w:t xml:space="preserve">or </w:t></w:r><w:r><w:t>linear (continuous)</w:t></w:r><w:r><w:t>. Hidden layers use ReLU, Leaky ReLU, Swish to prevent vanishing gradients in backprop. Input layers do not use activation functions.</w:t></w:r></w:p>'
Replace-ParagraphXml 43 $xml43

